# Update "paises" (countries) workbook per the 27-May-2020 07:35 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp string.
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 07:35"

# Pakistan (row 21): totals / new cases / active / recovered / deaths-today / deaths
$ws.Range("B21").Value = 59151
$ws.Range("C21").Value = 1446
$ws.Range("D21").Value = 19142
$ws.Range("E21").Value = 38784
$ws.Range("G21").Value = 28
$ws.Range("H21").Value = 1225

# Hungria (row 75)
$ws.Range("B75").Value = 3793
$ws.Range("C75").Value = 22
$ws.Range("D75").Value = 1856
$ws.Range("E75").Value = 1432
$ws.Range("G75").Value = 6
$ws.Range("H75").Value = 505

# Uzbekistan (row 76)
$ws.Range("B76").Value = 3333
$ws.Range("C76").Value = 43
$ws.Range("E76").Value = 683

# Rows 207/208 swap rank order: Islas Turcas y Caicos now outranks Groenlandia.
$ws.Range("A207").Value = "Islas Turcas y Caicos"
$ws.Range("D207").Value = 10
$ws.Range("H207").Value = 1

$ws.Range("A208").Value = "Groenlandia"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 0
